$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Whole Community" column headers (row 2) renamed to "Non-Verrucomicrobia"
$ws.Range("C2").Value = "Non-Verrucomicrobia"
$ws.Range("E2").Value = "Non-Verrucomicrobia"
$ws.Range("G2").Value = "Non-Verrucomicrobia"

# G2 picks up the same (bordered) formatting as the other header cells in row 2
$ws.Range("C2").Copy()
$ws.Range("G2").PasteSpecial(-4122)

# Updated "Best parameters(s)" values for the Estuary/Non-Verrucomicrobia columns
$ws.Range("C19").Value = "TDP"
$ws.Range("C20").Value = "TSS"
$ws.Range("G20").Value = "TDP"
$ws.Range("G21").ClearContents()
$ws.Range("G22").ClearContents()

# Updated correlation values in row 23
$ws.Range("C23").Value = "0.373 (0.02*)"
$ws.Range("E23").Value = "0.189 (0.39)"
$ws.Range("G23").Value = "0.063 (0.4)"

# Restore the active cell selection
$ws.Range("G24").Select()
